# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计",
#   populated with the Q1-2022 fund-holdings table (same shape/style as
#   the "2021-Q4" sheet).
# - Insert a new leading data row into "总计" for the "2022-Q1" summary
#   (date / holding count / holding value), pushing the existing
#   "2021-Q4" row down.
#
# NOTE: worksheet references here are re-fetched by name right before
# each use instead of being cached across structural edits (adding /
# moving sheets shifts the 1-based Index other variables may have
# resolved against).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted right after "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Match page margins / outline flags of the sibling sheets.
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1

# Copy the header-row / index-column look (bold, bordered, centered)
# from "2021-Q4" so the new sheet reuses the same cell style.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Item("2022-Q1")
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)

# Header row
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows: index column (numeric) + rank column (numeric) are plain
# numbers; everything else is the text exactly as recorded in the
# source table (codes/percentages kept as text so leading zeros and
# fixed decimals survive).
$rows = @(
    @(0, "000480", "东方红新动力灵活配置混合", "15.38", "72.90", "2.84", "0.4368", 9),
    @(1, "001564", "东方红京东大数据灵活配置混合", "11.40", "69.58", "2.68", "0.3055", 7),
    @(2, "006199", "长盛同锦研究精选混合", "1.73", "82.48", "2.87", "0.0497", 7),
    @(3, "001892", "长盛新兴成长主题灵活配置混合", "1.32", "82.10", "3.09", "0.0408", 6),
    @(4, "002085", "长盛互联网+主题灵活配置混合", "0.84", "83.97", "3.03", "0.0255", 6)
)

$q1 = $wb.Worksheets.Item("2022-Q1")
$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = "'" + $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Drop the quote-prefix styling picked up from the leading apostrophes
# above, leaving plain unstyled text cells (B:G) as in the source file.
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Range("B2:G6").ClearFormats()

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert the "2022-Q1" summary row above "2021-Q4"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total = $wb.Worksheets.Item("总计")
$total.Range("B2:D2").ClearFormats()
$total.Range("A2").Value = 0
$total.Range("B2").Value = "'2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.86
$total.Range("B2").ClearFormats()

# The pushed-down "2021-Q4" row's running index bumps from 0 to 1.
$total.Range("A3").Value = 1

# Restore the bold/bordered index-cell style on the new A2 (Insert only
# carried the border onto B2:D2, not A2).
$total = $wb.Worksheets.Item("总计")
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Leave the original sheet ("2021-Q4") as the active tab, as before.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Activate()
$q4.Range("A1").Select()
